$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.054.87"
$ws.Range("E2").Value = "  +0.91%  "

$ws.Range("D3").Value = "2.661.39"
$ws.Range("E3").Value = "  +4.93%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").Value = "'518.72"
$ws.Range("E5").Value = "  +2.38%  "

$ws.Range("D6").Value = "'144.96"
$ws.Range("E6").Value = "  +0.94%  "

$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("D8").Value = "'0.566"
$ws.Range("E8").Value = "  +0.80%  "

$ws.Range("D9").Value = "2.660.20"
$ws.Range("E9").Value = "  +4.63%  "

$ws.Range("D10").Value = "'6.28"
$ws.Range("E10").Value = "  +3.15%  "

$ws.Range("D11").Value = "'0.105"
$ws.Range("E11").Value = "  +2.85%  "

$ws.Range("D12").Value = "'0.337"
$ws.Range("E12").Value = "  +1.81%  "

$ws.Range("E13").Value = "  -1.69%  "

$ws.Range("D14").Value = "3.123.63"
$ws.Range("E14").Value = "  +4.71%  "

$ws.Range("D15").Value = "58.999.85"
$ws.Range("E15").Value = "  +0.81%  "

$ws.Range("D16").Value = "'20.97"
$ws.Range("E16").Value = "  +1.47%  "

$ws.Range("E17").Value = "  +1.76%  "

$ws.Range("D18").Value = "2.651.07"
$ws.Range("E18").Value = "  +4.33%  "

$ws.Range("D19").Value = "'351.42"
$ws.Range("E19").Value = "  +4.59%  "

$ws.Range("D20").Value = "'4.54"
$ws.Range("E20").Value = "  +0.17%  "

$ws.Range("D21").Value = "'10.36"
$ws.Range("E21").Value = "  +2.76%  "

$ws.Range("D22").Value = "'6.21"
$ws.Range("E22").Value = "  +4.45%  "

$ws.Range("E23").Value = "  -0.04%  "

$ws.Range("D24").Value = "'61.68"
$ws.Range("E24").Value = "  +2.00%  "

$ws.Range("D25").Value = "'0.420"
$ws.Range("E25").Value = "  +2.45%  "

$ws.Range("D26").Value = "2.752.33"
$ws.Range("E26").Value = "  +3.70%  "

$ws.Range("D27").Value = "'0.993"
$ws.Range("E27").Value = "  -0.56%  "

$ws.Range("E28").Value = "  +1.37%  "

$ws.Range("D29").Value = "0.0₃0806"
$ws.Range("E29").Value = "  +2.77%  "

$ws.Range("D30").Value = "'7.16"
$ws.Range("E30").Value = "  +3.17%  "

$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = "  +0.01%  "

$ws.Range("D32").Value = "'6.28"
$ws.Range("E32").Value = "  +7.86%  "

$ws.Range("D33").Value = "'19.02"
$ws.Range("E33").Value = "  +2.63%  "

$ws.Range("D34").Value = "'1.58"
$ws.Range("E34").Value = "  +2.96%  "

$ws.Range("E35").Value = "  -0.11%  "

$ws.Range("D36").Value = "'0.970"
$ws.Range("E36").Value = "  +6.25%  "

$ws.Range("D37").Value = "'4.02"
$ws.Range("E37").Value = "  +3.21%  "

$ws.Range("E38").Value = "  +2.46%  "

$ws.Range("D39").Value = "'36.72"
$ws.Range("E39").Value = "  +1.95%  "

$ws.Range("D40").Value = "'0.845"
$ws.Range("E40").Value = "  +3.07%  "

$ws.Range("D41").Value = "'3.72"
$ws.Range("E41").Value = "  +5.48%  "

$ws.Range("E42").Value = "  +1.77%  "

$ws.Range("D43").Value = "'278.73"
$ws.Range("E43").Value = "  -1.17%  "

$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "'0.998"
$ws.Range("E44").Value = "  -0.06%  "

$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").Value = "'0.612"
$ws.Range("E45").Value = "  +2.34%  "

$ws.Range("D46").Value = "'0.0984"
$ws.Range("E46").Value = "  -0.98%  "

$ws.Range("D47").Value = "'19.59"
$ws.Range("E47").Value = "  +5.22%  "

$ws.Range("D48").Value = "'0.0526"
$ws.Range("E48").Value = "  -1.09%  "

$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "'0.0230"
$ws.Range("E49").Value = "  +1.73%  "

$ws.Range("B50").Value = "WhiteBITCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D50").Value = "'10.30"
$ws.Range("E50").Value = "  +0.12%  "

$ws.Range("D51").Value = "1.990.97"
$ws.Range("E51").Value = "  +4.53%  "
